$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 'WATER BATH'
$ws.Range("B2").Value = 'ESSTELL'
$ws.Range("C2").Value = 'EWB-106D'
$ws.Range("D2").Value = 2018050863
$ws.Range("E2").Value = 'PTH01140'
$ws.Range("F2").Value = 45817
$ws.Range("G2").Value = 'WB25062805'

# --- Row 3 ---
$ws.Range("A3").Value = 'CENTRIFUGE'
$ws.Range("B3").Value = 'BIO-RAD'
$ws.Range("C3").Value = 'DIACENT-12'
$ws.Range("D3").Formula = '="2001047"'
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = 'PTH01142'
$ws.Range("F3").Value = 45817
$ws.Range("G3").Value = 'CF25062806'

# --- Row 4 ---
$ws.Range("A4").Value = 'CENTRIFUGE'
$ws.Range("B4").Value = 'CENTURION SCIENTIFIC'
$ws.Range("C4").Value = 'PRO-HOSPITAL.GP'
$ws.Range("D4").Value = '216069-12'
$ws.Range("E4").Value = 'PTH01143'
$ws.Range("F4").Value = 45817
$ws.Range("G4").Value = 'CF25062807'

# --- Row 5 ---
$ws.Range("A5").Value = 'CENTRIFUGE'
$ws.Range("B5").Value = 'NUVE'
$ws.Range("C5").Value = 'NF 200'
$ws.Range("D5").Formula = '="02.11738"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = 'PTH01144'
$ws.Range("F5").Value = 45817
$ws.Range("G5").Value = 'CF25062808'

# --- Row 6 ---
$ws.Range("A6").Value = 'HEMATOCRIT CENTRIFUGE'
$ws.Range("B6").Value = 'SCILOGEX'
$ws.Range("C6").Value = 'DM1424'
$ws.Range("D6").Value = 'LM4E000540'
$ws.Range("E6").Value = 'PTH01145'
$ws.Range("F6").Value = 45817
$ws.Range("G6").Value = 'CF25062809'

# --- Row 7 ---
$ws.Range("A7").Value = 'FREEZER'
$ws.Range("B7").Value = 'HAIER'
$ws.Range("C7").Value = 'HCF-300DP'
$ws.Range("D7").Value = 'B30LF 2E000 0QMM4 20019'
$ws.Range("E7").Value = 'PTH01146'
$ws.Range("F7").Value = 45817
$ws.Range("G7").Value = 'CH25062810'

# --- Row 8 ---
$ws.Range("A8").Value = 'DRY BATH INCUBATOR'
$ws.Range("B8").Value = 'MAJOR SCIENCE'
$ws.Range("C8").Value = 'EL-01'
$ws.Range("D8").Value = 130422246
$ws.Range("E8").Value = 'PTH01147'
$ws.Range("F8").Value = 45817
$ws.Range("G8").Value = 'DB25062811'

# --- Row 9 ---
$ws.Range("A9").Value = 'ROTATOR'
$ws.Range("B9").Value = 'DIGISYSTEM'
$ws.Range("C9").Value = 'DSR-2100A'
$ws.Range("D9").Formula = '="0910304"'
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = 'PTH01149'
$ws.Range("F9").Value = '-'
$ws.Range("G9").Value = 'ยกเลิก'

# --- Row 10 ---
$ws.Range("A10").Value = 'REFRIGERATOR'
$ws.Range("B10").Value = 'PHCBI'
$ws.Range("C10").Value = 'MBR-705GR-PE'
$ws.Range("D10").Formula = '="17100060"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = 'PTH01150'
$ws.Range("F10").Value = 45817
$ws.Range("G10").Value = 'CH25062813'

# --- Row 11 ---
$ws.Range("A11").Value = 'REFRIGERATOR'
$ws.Range("B11").Value = 'SONGSERM INTERCOOL'
$ws.Range("C11").Value = 'SDC-1000AY'
$ws.Range("D11").Value = 'SDC1000101-0212-03032'
$ws.Range("E11").Value = 'PTH01154'
$ws.Range("F11").Value = 45817
$ws.Range("G11").Value = 'CH25062814'

# --- Row 12 ---
$ws.Range("A12").Value = 'REFRIGERATOR'
$ws.Range("B12").Value = 'SIAMATIC'
$ws.Range("C12").Value = 'HURR PLUS 2 DS'
$ws.Range("D12").Value = 'HUP-02-L0166-1166-011'
$ws.Range("E12").Value = 'PTH01155'
$ws.Range("F12").Value = 45817
$ws.Range("G12").Value = 'CH25062815'

# --- Row 13 ---
$ws.Range("A13").Value = 'REFRIGERATOR'
$ws.Range("B13").Value = 'SIAMATIC'
$ws.Range("C13").Value = 'HURR PLUS 2 DS'
$ws.Range("D13").Value = 'HUP-02-L0166-1166-016'
$ws.Range("E13").Value = 'PTH01156'
$ws.Range("F13").Value = 45817
$ws.Range("G13").Value = 'CH25062816'

# --- Row 14 (previously blank placeholder row, now populated) ---
$ws.Range("A14").Value = 'ROTATOR'
$ws.Range("B14").Value = 'DLAB'
$ws.Range("C14").Value = 'SK-O180-S'
$ws.Range("D14").Value = 'SN24ABE0003461'
$ws.Range("E14").Value = '-'
$ws.Range("F14").Value = 45817
$ws.Range("G14").Value = 'CF25062812'

# --- View / selection changes ---
$excel.ActiveWindow.Zoom = 85
$ws.Range("R10").Select()
